$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 14773564.39128477
